# Fruta / hortaliza, semanal
#
# A new weekly price-report record needs to be inserted into the "Ajo"
# (garlic) price table. The new observation belongs right before the
# existing row 162, so every data row from 162 through 252 shifts down
# by one (becoming 163 through 253) and the sheet grows from
# A1:R252 to A1:R253.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push the existing row 162 (and everything below it) down by one row,
# opening up a blank row 162 for the new record.
$ws.Rows.Item(162).Insert()

# Populate the newly inserted row 162 with the new observation.
$ws.Cells.Item(162, 1).Value  = 9
$ws.Cells.Item(162, 2).Value  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(162, 3).Value  = "Metropolitana"
$ws.Cells.Item(162, 4).Value  = 44767
$ws.Cells.Item(162, 5).Value  = 13
$ws.Cells.Item(162, 6).Value  = 100112003
$ws.Cells.Item(162, 7).Value  = "Ajo"
$ws.Cells.Item(162, 8).Value  = "Chino"
$ws.Cells.Item(162, 9).Value  = "Primera"
$ws.Cells.Item(162, 10).Value = 520
$ws.Cells.Item(162, 11).Value = 27000
$ws.Cells.Item(162, 12).Value = 27000
$ws.Cells.Item(162, 13).Value = 27000
$ws.Cells.Item(162, 14).Value = "`$/caja 10 kilos"
$ws.Cells.Item(162, 15).Value = "China"
$ws.Cells.Item(162, 16).Value = 2700
$ws.Cells.Item(162, 17).Value = 10
$ws.Cells.Item(162, 18).Value = "Hortaliza"
